## jungnang-gu.xlsx — "Add files via upload"
##
## The uploaded workbook has one row removed from the apartment list:
## the entry for id 15902 ("세광") is deleted, and every row below it
## shifts up by one (the last row, id 824 / "한신", becomes the new
## last row, taking over row 74).
##
## Deleting the worksheet row reproduces this exactly:
##  - sheetData row 74 (15902/세광) is removed and old row 75 (824/한신)
##    becomes the new row 74
##  - the <dimension> shrinks from A1:B75 to A1:B74
##  - the now-unused "세광" shared string is dropped, shifting "한신"
##    into its slot so B74 still points at shared-string index 74

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(74).Delete()

# Match the saved cursor/selection position after the edit.
$ws.Range("H73").Select()
